$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16, shifting existing rows 16-45 down to 17-46
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with the new record's data.
# Columns A,B,C,E,F,G,H,I,Q,R carry the same values as the rest of this
# market/category block (identical to the old row 16, now row 17).
$ws.Cells.Item(16, 1).Value = 11
$ws.Cells.Item(16, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(16, 3).Value = "Bíobío"
$ws.Cells.Item(16, 4).Value = 44638
$ws.Cells.Item(16, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(16, 5).Value = 8
$ws.Cells.Item(16, 6).Value = 100112031
$ws.Cells.Item(16, 7).Value = "Poroto verde"
$ws.Cells.Item(16, 8).Value = "Magnum"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 180
$ws.Cells.Item(16, 11).Value = 23000
$ws.Cells.Item(16, 12).Value = 24000
$ws.Cells.Item(16, 13).Value = 23444
$ws.Cells.Item(16, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(16, 15).Value = "Región Metropolitana"
$ws.Cells.Item(16, 16).Value = 938
$ws.Cells.Item(16, 17).Value = 25
$ws.Cells.Item(16, 18).Value = "Hortaliza"
